$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

# ------------------------------------------------------------------
# The "Sample Type" and "Repeat Location" rows are being folded into
# the "Sampling Point ID" row (renamed "Sampling Point Type/Location").
# Delete them first (bottom-up so row indices for subsequent rows do
# not shift underneath us), then rewrite the surviving rows' text.
# ------------------------------------------------------------------
$t.Rows.Item(8).Delete()   # "Repeat Location"
$t.Rows.Item(7).Delete()   # "Sample Type"

# Row 6 is now "Sampling Point ID" - rename + rewrite description/notes.
$nameCell = $t.Cell(6, 1)
$nameCell.Range.Find.Execute("Sampling Point ID", $true, $false, $false, $false, $false, $true, 0, $false, "Sampling Point Type/Location", 1)

$descCell = $t.Cell(6, 2)
$descCell.Range.Find.Execute("Identifier for the sample station/location within the Water System Facility from which the sample is drawn.", $true, $false, $false, $false, $false, $true, 0, $false, "Indicate the type of sample. If a repeat sample, also indicate the repeat location.", 1)
$descCell2 = $t.Cell(6, 2)
$descCell2.Range.Find.Execute("(Required, 12 characters max.)", $true, $false, $false, $false, $false, $true, 0, $false, "(Required.)", 1)

$valCell = $t.Cell(6, 3)
$valCell.Range.Find.Execute("Must not be longer than 12 characters.", $true, $false, $false, $false, $false, $true, 0, $false, "Must select value from list.", 1)

# Rows 7 & 8 (previously 9 & 10) are "Original Lab Sample ID" and
# "Original Sample Collection Date" - unchanged apart from the
# "Sample Type = ""Repeat""" -> "Sample Type is Repeat" wording tweak
# in their "(Required if ...)" notes. Use wdReplaceOne (1, not the
# document-wide wdReplaceAll) so each Find stays scoped to its own
# cell and doesn't bleed into the other "(Required if ...)" cell.
$origIdReqCell = $t.Cell(7, 2)
$origIdReqCell.Range.Find.Execute('(Required if Sample Type = "Repeat". Otherwise, unused.)', $true, $false, $false, $false, $false, $true, 0, $false, "(Required if Sample Type is Repeat. Otherwise, unused.)", 1)

$origDateReqCell = $t.Cell(8, 2)
$origDateReqCell.Range.Find.Execute('(Required if Sample Type = "Repeat". Otherwise, unused.)', $true, $false, $false, $false, $false, $true, 0, $false, "(Required if Sample Type is Repeat. Otherwise, unused.)", 1)
